$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "RM 232" row (row 26) entirely.
$ws.Rows.Item(26).Delete()

# After the above deletion, the former "SC 92" row (originally row 28)
# has shifted up to row 27. Delete it entirely as well.
$ws.Rows.Item(27).Delete()

# Remaining rows have shifted up by two. Update the F-column (error) values
# that moved between the "SC 5", "SC 101" and "SC 119" rows.
$ws.Range("F26").Value = ""
$ws.Range("F27").Value = 17
$ws.Range("F29").Value = ""
